# Add payment 79174445 (Cash) 2025-08-23T09:41:43
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5's phone number (A5) was stored as text; copy it down to the new
# row 6 first so A6 keeps the original text representation ("79174445"),
# then convert A5 itself to a real number.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A5").Value = 79174445

# Fill in the rest of the new payment record on row 6.
$ws.Range("B6").Value = 5000
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 5000
$ws.Range("G6").Value = "Cash"
$ws.Range("H6").Value = "2025-08-23T09:41:43"
